$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Row 2
$ws.Range("B2").Value = "Bitcoin"
$ws.Range("C2").Value = "https://coinranking.com/coin/Qwsogvtv82FCd+bitcoin-btc"
$ws.Range("D2").Value = "25.866.21"
$ws.Range("E2").Value = "  -0.08%  "

# Row 3
$ws.Range("B3").Value = "Ethereum"
$ws.Range("C3").Value = "https://coinranking.com/coin/razxDUgYGNAdQ+ethereum-eth"
$ws.Range("D3").Value = "1.736.20"
$ws.Range("E3").Value = "  -0.29%  "

# Row 4
$ws.Range("B4").Value = "TetherUSD"
$ws.Range("C4").Value = "https://coinranking.com/coin/HIVsRcGKkPFtW+tetherusd-usdt"
$ws.Range("D4").NumberFormat = "@"
$ws.Range("D4").Value = "0.9997"
$ws.Range("E4").Value = "  -0.03%  "

# Row 5
$ws.Range("B5").Value = "BNB"
$ws.Range("C5").Value = "https://coinranking.com/coin/WcwrkfNI4FUAe+bnb-bnb"
$ws.Range("D5").NumberFormat = "@"
$ws.Range("D5").Value = "242.93"
$ws.Range("E5").Value = "  +5.19%  "

# Row 6
$ws.Range("B6").Value = "USDC"
$ws.Range("C6").Value = "https://coinranking.com/coin/aKzUVe4Hh_CON+usdc-usdc"
$ws.Range("D6").NumberFormat = "@"
$ws.Range("D6").Value = "1.000"
$ws.Range("E6").Value = "  +0.01%  "

# Row 7
$ws.Range("B7").Value = "XRP"
$ws.Range("C7").Value = "https://coinranking.com/coin/-l8Mn2pVlRs-p+xrp-xrp"
$ws.Range("D7").NumberFormat = "@"
$ws.Range("D7").Value = "0.5192"
$ws.Range("E7").Value = "  -1.16%  "

# Row 8
$ws.Range("B8").Value = "Cardano"
$ws.Range("C8").Value = "https://coinranking.com/coin/qzawljRxB5bYu+cardano-ada"
$ws.Range("D8").NumberFormat = "@"
$ws.Range("D8").Value = "0.2743"
$ws.Range("E8").Value = "  -0.28%  "

# Row 9
$ws.Range("B9").Value = "OKB"
$ws.Range("C9").Value = "https://coinranking.com/coin/PDKcptVnzJTmN+okb-okb"
$ws.Range("D9").NumberFormat = "@"
$ws.Range("D9").Value = "39.91"
$ws.Range("E9").Value = "  +1.33%  "

# Row 10
$ws.Range("B10").Value = "Dogecoin"
$ws.Range("C10").Value = "https://coinranking.com/coin/a91GCGd_u96cF+dogecoin-doge"
$ws.Range("D10").NumberFormat = "@"
$ws.Range("D10").Value = "0.06163"
$ws.Range("E10").Value = "  +0.33%  "

# Row 11
$ws.Range("B11").Value = "WrappedEther"
$ws.Range("C11").Value = "https://coinranking.com/coin/Mtfb0obXVh59u+wrappedether-weth"
$ws.Range("D11").Value = "1.743.98"
$ws.Range("E11").Value = "  +0.13%  "

# Row 12
$ws.Range("B12").Value = "TRON"
$ws.Range("C12").Value = "https://coinranking.com/coin/qUhEFk1I61atv+tron-trx"
$ws.Range("D12").NumberFormat = "@"
$ws.Range("D12").Value = "0.07194"
$ws.Range("E12").Value = "  +1.39%  "

# Row 13
$ws.Range("B13").Value = "Solana"
$ws.Range("C13").Value = "https://coinranking.com/coin/zNZHO_Sjf+solana-sol"
$ws.Range("D13").NumberFormat = "@"
$ws.Range("D13").Value = "14.97"
$ws.Range("E13").Value = "  -1.46%  "

# Row 14
$ws.Range("B14").Value = "Polygon"
$ws.Range("C14").Value = "https://coinranking.com/coin/uW2tk-ILY0ii+polygon-matic"
$ws.Range("D14").NumberFormat = "@"
$ws.Range("D14").Value = "0.6431"
$ws.Range("E14").Value = "  +0.15%  "

# Row 15
$ws.Range("B15").Value = "Polkadot"
$ws.Range("C15").Value = "https://coinranking.com/coin/25W7FG7om+polkadot-dot"
$ws.Range("D15").NumberFormat = "@"
$ws.Range("D15").Value = "4.616"
$ws.Range("E15").Value = "  +1.92%  "

# Row 16
$ws.Range("B16").Value = "Litecoin"
$ws.Range("C16").Value = "https://coinranking.com/coin/D7B1x_ks7WhV5+litecoin-ltc"
$ws.Range("D16").NumberFormat = "@"
$ws.Range("D16").Value = "77.31"
$ws.Range("E16").Value = "  -0.29%  "

# Row 17
$ws.Range("B17").Value = "Dai"
$ws.Range("C17").Value = "https://coinranking.com/coin/MoTuySvg7+dai-dai"
$ws.Range("D17").NumberFormat = "@"
$ws.Range("D17").Value = "1.000"
$ws.Range("E17").Value = "  +0.08%  "

# Row 18
$ws.Range("B18").Value = "BinanceUSD"
$ws.Range("C18").Value = "https://coinranking.com/coin/vSo2fu9iE1s0Y+binanceusd-busd"
$ws.Range("D18").NumberFormat = "@"
$ws.Range("D18").Value = "0.9995"
$ws.Range("E18").Value = "  -0.15%  "

# Row 19
$ws.Range("B19").Value = "WrappedBTC"
$ws.Range("C19").Value = "https://coinranking.com/coin/x4WXHge-vvFY+wrappedbtc-wbtc"
$ws.Range("D19").Value = "25.902.38"
$ws.Range("E19").Value = "  +0.08%  "

# Row 20
$ws.Range("B20").Value = "Avalanche"
$ws.Range("C20").Value = "https://coinranking.com/coin/dvUj0CzDZ+avalanche-avax"
$ws.Range("D20").NumberFormat = "@"
$ws.Range("D20").Value = "11.77"
$ws.Range("E20").Value = "  +2.22%  "

# Row 21
$ws.Range("B21").Value = "ShibaInu"
$ws.Range("C21").Value = "https://coinranking.com/coin/xz24e0BjL+shibainu-shib"
$ws.Range("D21").NumberFormat = "@"
$ws.Range("D21").Value = "0.000006795"
$ws.Range("E21").Value = "  +1.87%  "

# Row 22
$ws.Range("B22").Value = "WrappedliquidstakedEther2.0"
$ws.Range("C22").Value = "https://coinranking.com/coin/CiixT63n3+wrappedliquidstakedether20-wsteth"
$ws.Range("D22").Value = "1.962.61"
$ws.Range("E22").Value = "  +0.00%  "

# Row 23
$ws.Range("B23").Value = "Uniswap"
$ws.Range("C23").Value = "https://coinranking.com/coin/_H5FVG9iW+uniswap-uni"
$ws.Range("D23").NumberFormat = "@"
$ws.Range("D23").Value = "4.283"
$ws.Range("E23").Value = "  -0.24%  "

# Row 24
$ws.Range("B24").Value = "Cosmos"
$ws.Range("C24").Value = "https://coinranking.com/coin/Knsels4_Ol-Ny+cosmos-atom"
$ws.Range("D24").NumberFormat = "@"
$ws.Range("D24").Value = "8.640"
$ws.Range("E24").Value = "  -1.33%  "

# Row 25
$ws.Range("B25").Value = "Chainlink"
$ws.Range("C25").Value = "https://coinranking.com/coin/VLqpJwogdhHNb+chainlink-link"
$ws.Range("D25").NumberFormat = "@"
$ws.Range("D25").Value = "5.289"
$ws.Range("E25").Value = "  +2.57%  "

# Row 26
$ws.Range("B26").Value = "Monero"
$ws.Range("C26").Value = "https://coinranking.com/coin/3mVx2FX_iJFp5+monero-xmr"
$ws.Range("D26").NumberFormat = "@"
$ws.Range("D26").Value = "136.94"
$ws.Range("E26").Value = "  -2.45%  "

# Row 27
$ws.Range("B27").Value = "Toncoin"
$ws.Range("C27").Value = "https://coinranking.com/coin/67YlI0K1b+toncoin-ton"
$ws.Range("D27").NumberFormat = "@"
$ws.Range("D27").Value = "1.510"
$ws.Range("E27").Value = "  -0.67%  "

# Row 28
$ws.Range("B28").Value = "EthereumClassic"
$ws.Range("C28").Value = "https://coinranking.com/coin/hnfQfsYfeIGUQ+ethereumclassic-etc"
$ws.Range("D28").NumberFormat = "@"
$ws.Range("D28").Value = "15.22"
$ws.Range("E28").Value = "  +0.47%  "

# Row 29
$ws.Range("B29").Value = "LidoDAOToken"
$ws.Range("C29").Value = "https://coinranking.com/coin/Pe93bIOD2+lidodaotoken-ldo"
$ws.Range("D29").NumberFormat = "@"
$ws.Range("D29").Value = "1.772"
$ws.Range("E29").Value = "  -0.93%  "

# Row 30
$ws.Range("B30").Value = "BitcoinCash"
$ws.Range("C30").Value = "https://coinranking.com/coin/ZlZpzOJo43mIo+bitcoincash-bch"
$ws.Range("D30").NumberFormat = "@"
$ws.Range("D30").Value = "105.23"
$ws.Range("E30").Value = "  +2.46%  "

# Row 31
$ws.Range("B31").Value = "InternetComputer(DFINITY)"
$ws.Range("C31").Value = "https://coinranking.com/coin/aMNLwaUbY+internetcomputerdfinity-icp"
$ws.Range("D31").NumberFormat = "@"
$ws.Range("D31").Value = "3.960"
$ws.Range("E31").Value = "  +6.38%  "

# Row 32
$ws.Range("B32").Value = "Stellar"
$ws.Range("C32").Value = "https://coinranking.com/coin/f3iaFeCKEmkaZ+stellar-xlm"
$ws.Range("D32").NumberFormat = "@"
$ws.Range("D32").Value = "0.08245"
$ws.Range("E32").Value = "  -0.67%  "

# Row 33
$ws.Range("B33").Value = "Filecoin"
$ws.Range("C33").Value = "https://coinranking.com/coin/ymQub4fuB+filecoin-fil"
$ws.Range("D33").NumberFormat = "@"
$ws.Range("D33").Value = "3.645"
$ws.Range("E33").Value = "  +3.22%  "

# Row 34
$ws.Range("B34").Value = "Hedera"
$ws.Range("C34").Value = "https://coinranking.com/coin/jad286TjB+hedera-hbar"
$ws.Range("D34").NumberFormat = "@"
$ws.Range("D34").Value = "0.04682"
$ws.Range("E34").Value = "  +3.49%  "

# Row 35
$ws.Range("B35").Value = "HuobiToken"
$ws.Range("C35").Value = "https://coinranking.com/coin/DXwP4wF9ksbBO+huobitoken-ht"
$ws.Range("D35").NumberFormat = "@"
$ws.Range("D35").Value = "2.656"
$ws.Range("E35").Value = "  +1.62%  "

# Row 36
$ws.Range("B36").Value = "ARBITRUM"
$ws.Range("C36").Value = "https://coinranking.com/coin/1Uo6s62Oc+arbitrum-arb"
$ws.Range("D36").NumberFormat = "@"
$ws.Range("D36").Value = "0.9912"
$ws.Range("E36").Value = "  +1.46%  "

# Row 37
$ws.Range("B37").Value = "ImmutableX"
$ws.Range("C37").Value = "https://coinranking.com/coin/Z96jIvLU7+immutablex-imx"
$ws.Range("D37").NumberFormat = "@"
$ws.Range("D37").Value = "0.6199"
$ws.Range("E37").Value = "  -0.18%  "

# Row 38
$ws.Range("B38").Value = "MXToken"
$ws.Range("C38").Value = "https://coinranking.com/coin/QUC5kVAxSoB-+mxtoken-mx"
$ws.Range("D38").NumberFormat = "@"
$ws.Range("D38").Value = "2.688"
$ws.Range("E38").Value = "  +0.28%  "

# Row 39
$ws.Range("B39").Value = "VeChain"
$ws.Range("C39").Value = "https://coinranking.com/coin/FEbS54wxo4oIl+vechain-vet"
$ws.Range("D39").NumberFormat = "@"
$ws.Range("D39").Value = "0.01601"
$ws.Range("E39").Value = "  +0.60%  "

# Row 40
$ws.Range("B40").Value = "RenderToken"
$ws.Range("C40").Value = "https://coinranking.com/coin/7C4Mh4xy1yDel+rendertoken-rndr"
$ws.Range("D40").NumberFormat = "@"
$ws.Range("D40").Value = "1.924"
$ws.Range("E40").Value = "  -0.68%  "

# Row 41
$ws.Range("B41").Value = "PaxDollar"
$ws.Range("C41").Value = "https://coinranking.com/coin/JCKLgWPAF+paxdollar-usdp"
$ws.Range("D41").NumberFormat = "@"
$ws.Range("D41").Value = "0.9997"
$ws.Range("E41").Value = "  +0.03%  "

# Row 42
$ws.Range("B42").Value = "Quant"
$ws.Range("C42").Value = "https://coinranking.com/coin/bauj_21eYVwso+quant-qnt"
$ws.Range("D42").NumberFormat = "@"
$ws.Range("D42").Value = "99.83"
$ws.Range("E42").Value = "  -0.22%  "

# Row 43
$ws.Range("B43").Value = "TheSandbox"
$ws.Range("C43").Value = "https://coinranking.com/coin/pxtKbG5rg+thesandbox-sand"
$ws.Range("D43").NumberFormat = "@"
$ws.Range("D43").Value = "0.3853"
$ws.Range("E43").Value = "  -0.37%  "

# Row 44
$ws.Range("B44").Value = "TrustWalletToken"
$ws.Range("C44").Value = "https://coinranking.com/coin/Hm3OlynlC+trustwallettoken-twt"
$ws.Range("D44").NumberFormat = "@"
$ws.Range("D44").Value = "0.7468"
$ws.Range("E44").Value = "  +2.18%  "

# Row 45
$ws.Range("B45").Value = "FraxShare"
$ws.Range("C45").Value = "https://coinranking.com/coin/3nNpuxHJ8+fraxshare-fxs"
$ws.Range("D45").NumberFormat = "@"
$ws.Range("D45").Value = "5.011"
$ws.Range("E45").Value = "  +0.09%  "

# Row 46
$ws.Range("B46").Value = "Algorand"
$ws.Range("C46").Value = "https://coinranking.com/coin/TpHE2IShQw-sJ+algorand-algo"
$ws.Range("D46").NumberFormat = "@"
$ws.Range("D46").Value = "0.1128"
$ws.Range("E46").Value = "  +0.16%  "

# Row 47
$ws.Range("B47").Value = "Aptos"
$ws.Range("C47").Value = "https://coinranking.com/coin/HGYj5JCv5+aptos-apt"
$ws.Range("D47").NumberFormat = "@"
$ws.Range("D47").Value = "6.250"
$ws.Range("E47").Value = "  -0.08%  "

# Row 48
$ws.Range("B48").Value = "Aave"
$ws.Range("C48").Value = "https://coinranking.com/coin/ixgUfzmLR+aave-aave"
$ws.Range("D48").NumberFormat = "@"
$ws.Range("D48").Value = "55.04"
$ws.Range("E48").Value = "  +2.54%  "

# Row 49
$ws.Range("B49").Value = "Cronos"
$ws.Range("C49").Value = "https://coinranking.com/coin/65PHZTpmE55b+cronos-cro"
$ws.Range("D49").NumberFormat = "@"
$ws.Range("D49").Value = "0.05214"
$ws.Range("E49").Value = "  -2.31%  "

# Row 50
$ws.Range("B50").Value = "Elrond"
$ws.Range("C50").Value = "https://coinranking.com/coin/omwkOTglq+elrond-egld"
$ws.Range("D50").NumberFormat = "@"
$ws.Range("D50").Value = "30.63"
$ws.Range("E50").Value = "  +1.72%  "

# Row 51
$ws.Range("B51").Value = "EnergySwap"
$ws.Range("C51").Value = "https://coinranking.com/coin/SbWqqTui-+energyswap-ens"
$ws.Range("D51").NumberFormat = "@"
$ws.Range("D51").Value = "7.562"
$ws.Range("E51").Value = "  -1.13%  "
